# Update latest output (run 189)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Schedule sheet - summary cost figures
$wsSchedule.Range("E2").Value = 1196.95129125
$wsSchedule.Range("F2").Value = 19.79086129712301

# Detailed sheet - Price column (B) updates, plus a couple of Type (C) updates
$wsDetailed.Range("B17").Value = 34.69723

$wsDetailed.Range("B18").Value = 25.178

$wsDetailed.Range("B19").Value = 36.06
$wsDetailed.Range("C19").Value = "historical"

$wsDetailed.Range("B20").Value = 31.28559
$wsDetailed.Range("C20").Value = "historical"

$wsDetailed.Range("B21").Value = 34.35704

$wsDetailed.Range("B22").Value = 30.95505

$wsDetailed.Range("B23").Value = 31.61438

$wsDetailed.Range("B24").Value = 34.01

$wsDetailed.Range("B25").Value = 23.10783

$wsDetailed.Range("B26").Value = 21.66239

$wsDetailed.Range("B27").Value = 24.20626

$wsDetailed.Range("B28").Value = 31.16386

$wsDetailed.Range("B29").Value = 26.63675

$wsDetailed.Range("B30").Value = 23.41793

$wsDetailed.Range("B31").Value = 4.2416

$wsDetailed.Range("B32").Value = -7.37769

$wsDetailed.Range("B33").Value = -5.16197

$wsDetailed.Range("B34").Value = 0

$wsDetailed.Range("B35").Value = -4.3238

$wsDetailed.Range("B36").Value = 0.51

$wsDetailed.Range("B37").Value = 35.26909

$wsDetailed.Range("B38").Value = 57.03052

$wsDetailed.Range("B39").Value = 57.05891

$wsDetailed.Range("B40").Value = 57.09

$wsDetailed.Range("B41").Value = 58.81437

$wsDetailed.Range("B42").Value = 59.32526

$wsDetailed.Range("B44").Value = 57.06

$wsDetailed.Range("B46").Value = 57.06
